$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows with payment data added (F/G), while A/D/E already existed ---
$ws.Range("F20").Value = 44733
$ws.Range("G20").Value = 251

$ws.Range("F24").Value = 44699
$ws.Range("G24").Value = 3400

$ws.Range("F25").Value = 44726
$ws.Range("G25").Value = 10409

$ws.Range("F26").Value = 44733
$ws.Range("G26").Value = 1093

# --- Rows 27-43: full new credit entries (A, D, E, F, G) ---
$ws.Range("A27").Value = 44726
$ws.Range("D27").Value = "HERRADURA DAVID"
$ws.Range("E27").Value = 8615
$ws.Range("F27").Value = 44728
$ws.Range("G27").Value = 8615

$ws.Range("A28").Value = 44726
$ws.Range("D28").Value = "MARCELO"
$ws.Range("E28").Value = 6901
$ws.Range("F28").Value = 44727
$ws.Range("G28").Value = 6901

$ws.Range("A29").Value = 44726
$ws.Range("D29").Value = "HERRADURA GUSTAVO"
$ws.Range("E29").Value = 16303
$ws.Range("F29").Value = 44728
$ws.Range("G29").Value = 16303

$ws.Range("A30").Value = 44728
$ws.Range("D30").Value = "HERRADURA GUSTAVO"
$ws.Range("E30").Value = 13075
$ws.Range("F30").Value = 44730
$ws.Range("G30").Value = 13075

$ws.Range("A31").Value = 44728
$ws.Range("D31").Value = "HERRADURA DAVID"
$ws.Range("E31").Value = 11575
$ws.Range("F31").Value = 44729
$ws.Range("G31").Value = 11575

$ws.Range("A32").Value = 44729
$ws.Range("D32").Value = "HERRADURA DAVID"
$ws.Range("E32").Value = 18462
$ws.Range("F32").Value = 44730
$ws.Range("G32").Value = 18462

$ws.Range("A33").Value = 44729
$ws.Range("D33").Value = "HERRADURA DAVID"
$ws.Range("E33").Value = 5618
$ws.Range("F33").Value = 44730
$ws.Range("G33").Value = 5618

$ws.Range("A34").Value = 44730
$ws.Range("D34").Value = "OBRADOR"
$ws.Range("E34").Value = 329
$ws.Range("F34").Value = 44733
$ws.Range("G34").Value = 329

$ws.Range("A35").Value = 44730
$ws.Range("D35").Value = "HERRADURA DAVID"
$ws.Range("E35").Value = 20983
$ws.Range("F35").Value = 44731
$ws.Range("G35").Value = 20983

$ws.Range("A36").Value = 44730
$ws.Range("D36").Value = "HERRADURA GUSTAVO"
$ws.Range("E36").Value = 14363
$ws.Range("F36").Value = 44732
$ws.Range("G36").Value = 14363

$ws.Range("A37").Value = 44730
$ws.Range("D37").Value = "HERRADURA DAVID"
$ws.Range("E37").Value = 2957
$ws.Range("F37").Value = 44730
$ws.Range("G37").Value = 2957

$ws.Range("A38").Value = 44731
$ws.Range("D38").Value = "HERRADURA DAVID"
$ws.Range("E38").Value = 12874
$ws.Range("F38").Value = 44733
$ws.Range("G38").Value = 12874
# F38/G38 originally used the row-block's blue-font style (s35/s36); normalize
# them to the common style (s24/s25) used by every other F/G cell, matching
# row 37's formatting.
$ws.Range("F37:G37").Copy()
$ws.Range("F38:G38").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A39").Value = 44732
$ws.Range("D39").Value = "HERRADURA DAVID"
$ws.Range("E39").Value = 13802
$ws.Range("F39").Value = 44733
$ws.Range("G39").Value = 13802

$ws.Range("A40").Value = 44732
$ws.Range("D40").Value = "HERRADURA GUSTAVO"
$ws.Range("E40").Value = 16448
$ws.Range("F40").Value = 44735
$ws.Range("G40").Value = 16448

$ws.Range("A41").Value = 44733
$ws.Range("D41").Value = "HERRADURA DAVID"
$ws.Range("E41").Value = 11463
$ws.Range("F41").Value = 44736
$ws.Range("G41").Value = 11463

$ws.Range("A42").Value = 44733
$ws.Range("D42").Value = "MARCELO"
$ws.Range("E42").Value = 7538
$ws.Range("F42").Value = 44734
$ws.Range("G42").Value = 7538

$ws.Range("A43").Value = 44734
$ws.Range("D43").Value = "MARCELO"
$ws.Range("E43").Value = 3409
$ws.Range("F43").Value = 44736
$ws.Range("G43").Value = 3409

# --- Rows 44-46: new credit entries with no payment yet (F/G stay blank) ---
$ws.Range("A44").Value = 44735
$ws.Range("D44").Value = "OBRADOR"
$ws.Range("E44").Value = 268

$ws.Range("A45").Value = 44736
$ws.Range("D45").Value = "HERRADURA DAVID"
$ws.Range("E45").Value = 16050

$ws.Range("A46").Value = 44736
$ws.Range("D46").Value = "MAURO"
$ws.Range("E46").Value = 3129

# --- View state: scroll position + selection ---
$excel.Goto($ws.Range("A27"), $true)
$ws.Range("D47").Select()
